$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 363.81818
$ws.Range("I5").Value = 411.5
$ws.Range("J5").Value = 236.66667
$ws.Range("K5").Value = 411.5
$ws.Range("L5").Value = 236.66667
$ws.Range("M5").Value = -296.5
$ws.Range("N5").Value = -466.66667
$ws.Range("H18").Value = 1050
$ws.Range("I18").Value = 533.3333
$ws.Range("J18").Value = 1566.6666
$ws.Range("K18").Value = 533.3333
$ws.Range("L18").Value = 1566.6666
$ws.Range("M18").Value = -249.3333
$ws.Range("N18").Value = -2134.6666
$ws.Range("H33").Value = 3683.6206
$ws.Range("I33").Value = 4213
$ws.Range("J33").Value = 375
$ws.Range("K33").Value = 4213
$ws.Range("L33").Value = 375
$ws.Range("M33").Value = -3984
$ws.Range("N33").Value = -833
$ws.Range("H41").Value = 378.13333
$ws.Range("I41").Value = 459.81818
$ws.Range("J41").Value = 153.5
$ws.Range("K41").Value = 459.81818
$ws.Range("L41").Value = 153.5
$ws.Range("M41").Value = -19.81817999999998
$ws.Range("N41").Value = -1033.5
$ws.Range("H113").Value = 5424.7334
$ws.Range("I113").Value = 3852.5
$ws.Range("J113").Value = 5666.615
$ws.Range("K113").Value = 3852.5
$ws.Range("L113").Value = 5666.615
$ws.Range("M113").Value = -598.5
$ws.Range("N113").Value = -12174.615
$ws.Range("H129").Value = 1053.0944
$ws.Range("I129").Value = 506.2
$ws.Range("J129").Value = 1110.0625
$ws.Range("K129").Value = 1518.6
$ws.Range("L129").Value = 3330.1875
$ws.Range("M129").Value = 3481.4
$ws.Range("N129").Value = -13330.1875
$ws.Range("H132").Value = 27560.104
$ws.Range("I132").Value = 28928.217
$ws.Range("J132").Value = 2250
$ws.Range("K132").Value = 86784.651
$ws.Range("L132").Value = 6750
$ws.Range("M132").Value = -84254.651
$ws.Range("N132").Value = -11810
$ws.Range("H138").Value = 2636.6667
$ws.Range("I138").Value = 1290.6346
$ws.Range("J138").Value = 4636.486
$ws.Range("K138").Value = 3871.9038
$ws.Range("L138").Value = 13909.458
$ws.Range("M138").Value = 1268.0962
$ws.Range("N138").Value = -24189.458
$ws.Range("H141").Value = 6106.979
$ws.Range("I141").Value = 1261.0286
$ws.Range("J141").Value = 19153.77
$ws.Range("K141").Value = 3783.0858
$ws.Range("L141").Value = 57461.31
$ws.Range("M141").Value = 1396.9142
$ws.Range("N141").Value = -67821.31

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H18").Value = 50000
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 50000
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 50000
$ws.Range("M18").ClearContents()
$ws.Range("N18").Value = -50644
$ws.Range("H32").Value = 2089.08
$ws.Range("I32").Value = 1991.6907
$ws.Range("J32").Value = 5238
$ws.Range("K32").Value = 1991.6907
$ws.Range("L32").Value = 5238
$ws.Range("M32").Value = -1704.6907
$ws.Range("N32").Value = -5812
$ws.Range("H35").Value = 1000
$ws.Range("I35").Value = 1000
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 1000
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -594
$ws.Range("H97").Value = 1230.3077
$ws.Range("I97").Value = 922.2381
$ws.Range("J97").Value = 2524.2
$ws.Range("K97").Value = 922.2381
$ws.Range("L97").Value = 2524.2
$ws.Range("M97").Value = -426.2381
$ws.Range("N97").Value = -3516.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 12745.167
$ws.Range("I26").Value = 12745.167
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 12745.167
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -12453.167
$ws.Range("H99").Value = 2212.5
$ws.Range("I99").Value = 1750
$ws.Range("J99").Value = 2590.9092
$ws.Range("K99").Value = 1750
$ws.Range("L99").Value = 2590.9092
$ws.Range("M99").Value = -252
$ws.Range("N99").Value = -5586.9092

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 248.4
$ws.Range("I22").Value = 97.333336
$ws.Range("J22").Value = 475
$ws.Range("K22").Value = 97.333336
$ws.Range("L22").Value = 475
$ws.Range("M22").Value = 252.666664
$ws.Range("N22").Value = -1175
$ws.Range("H26").Value = 10000
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 10000
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 10000
$ws.Range("N26").Value = -10574
$ws.Range("H31").Value = 2499.5
$ws.Range("I31").Value = 1661.7931
$ws.Range("J31").Value = 3283.1614
$ws.Range("K31").Value = 1661.7931
$ws.Range("L31").Value = 3283.1614
$ws.Range("M31").Value = -1366.7931
$ws.Range("N31").Value = -3873.1614
$ws.Range("H34").Value = 2499.5
$ws.Range("I34").Value = 1661.7931
$ws.Range("J34").Value = 3283.1614
$ws.Range("K34").Value = 1661.7931
$ws.Range("L34").Value = 3283.1614
$ws.Range("M34").Value = -1459.7931
$ws.Range("N34").Value = -3687.1614
$ws.Range("H56").Value = 50000
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 50000
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 50000
$ws.Range("N56").Value = -51690
$ws.Range("H58").Value = 1404.6072
$ws.Range("I58").Value = 1103.4131
$ws.Range("J58").Value = 2790.1
$ws.Range("K58").Value = 1103.4131
$ws.Range("L58").Value = 2790.1
$ws.Range("M58").Value = -900.4131
$ws.Range("N58").Value = -3196.1
$ws.Range("H136").Value = 1404.6072
$ws.Range("I136").Value = 1103.4131
$ws.Range("J136").Value = 2790.1
$ws.Range("K136").Value = 3310.2393
$ws.Range("L136").Value = 8370.299999999999
$ws.Range("M136").Value = -760.2393000000002
$ws.Range("N136").Value = -13470.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 975.1667
$ws.Range("I34").Value = 498
$ws.Range("J34").Value = 1070.6
$ws.Range("K34").Value = 1494
$ws.Range("L34").Value = 3211.8
$ws.Range("M34").Value = -1410
$ws.Range("N34").Value = -3379.8
$ws.Range("H113").Value = 1242.7858
$ws.Range("I113").Value = 3126.75
$ws.Range("J113").Value = 489.2
$ws.Range("K113").Value = 9380.25
$ws.Range("L113").Value = 1467.6
$ws.Range("M113").Value = -7210.25
$ws.Range("N113").Value = -5807.6
$ws.Range("H131").Value = 3518.805
$ws.Range("I131").Value = 303.44446
$ws.Range("J131").Value = 4423.125
$ws.Range("K131").Value = 910.33338
$ws.Range("L131").Value = 13269.375
$ws.Range("M131").Value = 4129.66662
$ws.Range("N131").Value = -23349.375
$ws.Range("H132").Value = 1807.8422
$ws.Range("I132").Value = 2011.7778
$ws.Range("J132").Value = 1624.3
$ws.Range("K132").Value = 18106.0002
$ws.Range("L132").Value = 14618.7
$ws.Range("M132").Value = -15576.0002
$ws.Range("N132").Value = -19678.7

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 1031
$ws.Range("I31").Value = 1031
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 1031
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -739
$ws.Range("H37").Value = 1031
$ws.Range("I37").Value = 1031
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 1031
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -754
$ws.Range("H124").Value = 39200
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 39200
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 39200
$ws.Range("N124").Value = -49020
$ws.Range("H126").Value = 1634.7142
$ws.Range("I126").Value = 1472.9
$ws.Range("J126").Value = 2039.25
$ws.Range("K126").Value = 4418.700000000001
$ws.Range("L126").Value = 6117.75
$ws.Range("M126").Value = -1948.700000000001
$ws.Range("N126").Value = -11057.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 20000
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 20000
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 20000
$ws.Range("N43").Value = -20386
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("N45").ClearContents()
$ws.Range("H132").Value = 5374.9624
$ws.Range("I132").Value = 5686.18
$ws.Range("J132").Value = 4856.2666
$ws.Range("K132").Value = 17058.54
$ws.Range("L132").Value = 14568.7998
$ws.Range("M132").Value = -14528.54
$ws.Range("N132").Value = -19628.7998
$ws.Range("H136").Value = 3286.75
$ws.Range("I136").Value = 3613.0715
$ws.Range("J136").Value = 1002.5
$ws.Range("K136").Value = 10839.2145
$ws.Range("L136").Value = 3007.5
$ws.Range("M136").Value = -8289.2145
$ws.Range("N136").Value = -8107.5
$ws.Range("H141").Value = 53360.625
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 53360.625
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 53360.625
$ws.Range("N141").Value = -63720.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1137.44
$ws.Range("I136").Value = 567.2162
$ws.Range("J136").Value = 2760.3845
$ws.Range("K136").Value = 1701.6486
$ws.Range("L136").Value = 8281.1535
$ws.Range("M136").Value = 848.3514
$ws.Range("N136").Value = -13381.1535
